# [Kadastro App] Yeni kayit eklendi: 2938
#
# Appends the new "Erdemli" record (Kayit No 2938) to the master
# "Kayitlar" log sheet and to the per-birim "Erdemli" sheet, mirroring
# how every other record already appears in both places.

$wb = $excel.ActiveWorkbook

$kayitNo   = "2938"
$tarih     = "2025-09-09"
$birim     = "Erdemli"
$parsel    = "5"
$is        = "ÇAP"
$personel  = "SEVİL SARAÇER (Tekniker)"

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the next empty row right after the existing records.
    $lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
    $newRow = $lastRow + 1

    # Force text storage for values that otherwise look like a number or a
    # date (matches the other rows already in the sheet, which are all
    # stored as plain text), plain words don't need the quote-prefix.
    $ws.Cells.Item($newRow, 1).Value = "'" + $kayitNo
    $ws.Cells.Item($newRow, 2).Value = "'" + $tarih
    $ws.Cells.Item($newRow, 3).Value = $birim
    $ws.Cells.Item($newRow, 4).Value = "'" + $parsel
    $ws.Cells.Item($newRow, 5).Value = $is
    $ws.Cells.Item($newRow, 6).Value = $personel
}
